# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.595.76"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "2.614.81"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.582"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.72%  "
$ws.Range("D9").Value = "2.617.28"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.89%  "
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.333"
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "3.070.34"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "58.542.85"
$ws.Range("E15").Value = "  -1.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "2.613.29"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "334.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.421"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -2.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("D29").Value = "0.0₃0735"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "37.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.833"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.823"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.56%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "283.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.596"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0945"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0529"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0226"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.939.88"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.75%  "
